$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

Set-TextValue "D2" "29.417.85"
Set-TextValue "E2" "  +0.42%  "
Set-TextValue "D3" "1.876.74"
Set-TextValue "E3" "  +0.85%  "
Set-TextValue "E4" "  +0.19%  "
Set-TextValue "D5" "0.7173"
Set-TextValue "E5" "  +0.81%  "
Set-TextValue "D6" "240.10"
Set-TextValue "E6" "  +0.93%  "
Set-TextValue "D7" "0.9998"
Set-TextValue "E7" "  +0.14%  "
Set-TextValue "D8" "0.07838"
Set-TextValue "E8" "  -1.61%  "
Set-TextValue "D9" "0.3098"
Set-TextValue "E9" "  +1.86%  "
Set-TextValue "D10" "24.80"
Set-TextValue "E10" "  +5.07%  "
Set-TextValue "D11" "0.08256"
Set-TextValue "E11" "  +1.01%  "
Set-TextValue "B12" "Polygon"
Set-TextValue "C12" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D12" "0.7267"
Set-TextValue "E12" "  +3.04%  "
Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "5.283"
Set-TextValue "E13" "  +1.90%  "
Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.858.06"
Set-TextValue "E14" "  +0.73%  "
Set-TextValue "D15" "91.24"
Set-TextValue "E15" "  +1.60%  "
Set-TextValue "D16" "29.410.86"
Set-TextValue "E16" "  +0.51%  "
Set-TextValue "D17" "5.911"
Set-TextValue "E17" "  +0.94%  "
Set-TextValue "D18" "245.15"
Set-TextValue "E18" "  +3.10%  "
Set-TextValue "E19" "  -0.03%  "
Set-TextValue "D20" "13.30"
Set-TextValue "E20" "  -0.15%  "
Set-TextValue "E21" "  +0.30%  "
Set-TextValue "D22" "7.901"
Set-TextValue "E22" "  +6.54%  "
Set-TextValue "D23" "0.9998"
Set-TextValue "E23" "  +0.23%  "
Set-TextValue "D24" "0.1564"
Set-TextValue "E24" "  +8.52%  "
Set-TextValue "D25" "163.93"
Set-TextValue "E25" "  +1.10%  "
Set-TextValue "D26" "9.016"
Set-TextValue "E26" "  +0.69%  "
Set-TextValue "D27" "18.35"
Set-TextValue "E27" "  +1.37%  "
Set-TextValue "D28" "1.362"
Set-TextValue "E28" "  -4.97%  "
Set-TextValue "D29" "1.484"
Set-TextValue "E29" "  +0.33%  "
Set-TextValue "D30" "4.396"
Set-TextValue "E30" "  +0.25%  "
Set-TextValue "D31" "4.144"
Set-TextValue "E31" "  +3.13%  "
Set-TextValue "D32" "0.05281"
Set-TextValue "E32" "  +1.11%  "
Set-TextValue "D33" "1.935"
Set-TextValue "E33" "  +0.10%  "
Set-TextValue "D34" "1.201"
Set-TextValue "E34" "  +3.09%  "
Set-TextValue "D35" "0.7224"
Set-TextValue "E35" "  +1.65%  "
Set-TextValue "D36" "2.677"
Set-TextValue "E36" "  +0.58%  "
Set-TextValue "D37" "0.01868"
Set-TextValue "E37" "  +0.95%  "
Set-TextValue "D38" "1.232.47"
Set-TextValue "E38" "  +9.38%  "
Set-TextValue "E39" "  -0.19%  "
Set-TextValue "D40" "0.9094"
Set-TextValue "E40" "  -2.18%  "
Set-TextValue "D41" "73.01"
Set-TextValue "E41" "  +4.40%  "
Set-TextValue "D42" "6.077"
Set-TextValue "E42" "  +3.89%  "
Set-TextValue "D43" "0.9998"
Set-TextValue "E43" "  +0.15%  "
Set-TextValue "D44" "103.74"
Set-TextValue "E44" "  +0.55%  "
Set-TextValue "D45" "0.5339"
Set-TextValue "E45" "  -0.11%  "
Set-TextValue "B46" "RocketPoolETH"
Set-TextValue "C46" "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue "D46" "2.015.18"
Set-TextValue "E46" "  +1.13%  "
Set-TextValue "D47" "2.943"
Set-TextValue "E47" "  +12.48%  "
Set-TextValue "B48" "BabyDogeCoin"
Set-TextValue "C48" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D48" "0.00000000121"
Set-TextValue "E48" "  +2.23%  "
Set-TextValue "B49" "RenderToken"
Set-TextValue "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D49" "1.758"
Set-TextValue "E49" "  -0.44%  "
Set-TextValue "B50" "TheSandbox"
Set-TextValue "C50" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D50" "0.4328"
Set-TextValue "E50" "  +1.22%  "
Set-TextValue "B51" "EnergySwap"
Set-TextValue "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "9.255"

Write-Host "Applied all cell updates"
